# "add tieumuc vao muchi" - add a "Tiểu mục" column into the mục chi (expense item) table.
#
# The original sheet1 layout (row 2 = header row) is:
#   A: Mã nhóm hoạt động | B: Tên | C: Tên mục chi | D: Ghi chú
#
# The edit inserts a brand-new column before column D ("Ghi chú"), which shifts
# "Ghi chú" from D to E, and fills the new D column's header (row 2) with the
# new shared string "Tiểu mục". The header row/title merge (A1:D1) grows to
# A1:E1, the sheet dimension grows from A1:D4 to A1:E4, and the active
# selection ends up on D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at D; this shifts the old "Ghi chú" column (D) to E
# and extends row1's merged title cell / sheet dimension automatically.
$ws.Columns("D").Insert() | Out-Null

# The two data rows (3 and 4) only had data through column C; after the
# column insert they pick up an empty placeholder cell in the new column D.
# Remove it so those rows stay exactly as they were (no stray D cell).
$ws.Range("D3:D4").Clear() | Out-Null

# Give the new column a header: "Tiểu mục".
$ws.Range("D2").Value = "Tiểu mục"

# Match the final selection left behind in the saved workbook.
$ws.Range("D9").Select() | Out-Null
